$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link URL) - plain string, no numeric coercion risk ---
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"

# --- Numeric-looking text columns (Price / Volume%) ---
# Force text storage (NumberFormat "@") so literal formatting (trailing zeros, "%", "," etc.)
# survives instead of being coerced into a number, then clear the format so the cell
# style matches the original (no explicit style index).
$numericTextCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "E47", "D48", "E48", "E49", "E50")
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "254.88"
$ws.Range("E2").Value = "3.67%"
$ws.Range("D3").Value = "28.09"
$ws.Range("E3").Value = "-5.59%"
$ws.Range("D4").Value = "5.292"
$ws.Range("E4").Value = "2.70%"
$ws.Range("D5").Value = "0.05848"
$ws.Range("E5").Value = "1.32%"
$ws.Range("E6").Value = "0.61%"
$ws.Range("D7").Value = "0.8700"
$ws.Range("E7").Value = "1.93%"
$ws.Range("D8").Value = "0.9371"
$ws.Range("E8").Value = "9.77%"
$ws.Range("D9").Value = "0.1412"
$ws.Range("E9").Value = "2.51%"
$ws.Range("D10").Value = "0.07102"
$ws.Range("E10").Value = "0.21%"
$ws.Range("D11").Value = "0.03177"
$ws.Range("E11").Value = "-2.62%"
$ws.Range("D12").Value = "0.09225"
$ws.Range("E12").Value = "-1.50%"
$ws.Range("D13").Value = "0.001551"
$ws.Range("E13").Value = "1.32%"
$ws.Range("D14").Value = "0.005802"
$ws.Range("E14").Value = "-3.87%"
$ws.Range("D15").Value = "3.498"
$ws.Range("E15").Value = "-0.38%"
$ws.Range("D16").Value = "3.231"
$ws.Range("E16").Value = "-0.24%"
$ws.Range("D17").Value = "2.222"
$ws.Range("E17").Value = "-0.02%"
$ws.Range("D18").Value = "0.01062"
$ws.Range("E18").Value = "1,680.86%"
$ws.Range("D19").Value = "0.3180"
$ws.Range("E19").Value = "0.70%"
$ws.Range("D20").Value = "0.03460"
$ws.Range("E20").Value = "2.51%"
$ws.Range("E21").Value = "1.47%"
$ws.Range("E22").Value = "0.95%"
$ws.Range("D23").Value = "0.04155"
$ws.Range("E23").Value = "0.76%"
$ws.Range("D24").Value = "0.1345"
$ws.Range("E24").Value = "-4.62%"
$ws.Range("D25").Value = "0.001229"
$ws.Range("E25").Value = "0.01%"
$ws.Range("D26").Value = "0.005002"
$ws.Range("E26").Value = "20.85%"
$ws.Range("D27").Value = "0.0001200"
$ws.Range("E27").Value = "-0.06%"
$ws.Range("D28").Value = "0.00007999"
$ws.Range("E28").Value = "-44.82%"
$ws.Range("E40").Value = "1.84%"
$ws.Range("D41").Value = "0.1100"
$ws.Range("E41").Value = "2.86%"
$ws.Range("D42").Value = "0.003813"
$ws.Range("E42").Value = "-32.66%"
$ws.Range("D43").Value = "0.002300"
$ws.Range("E43").Value = "-0.06%"
$ws.Range("D44").Value = "0.01008"
$ws.Range("E44").Value = "18.69%"
$ws.Range("D45").Value = "0.00005225"
$ws.Range("E45").Value = "-3.78%"
$ws.Range("E46").Value = "-0.06%"
$ws.Range("D47").Value = "0.09299"
$ws.Range("E47").Value = "30.91%"
$ws.Range("D48").Value = "0.002154"
$ws.Range("E48").Value = "-2.74%"
$ws.Range("E49").Value = "-0.06%"
$ws.Range("E50").Value = "-0.06%"

foreach ($ref in $numericTextCells) {
    $ws.Range($ref).ClearFormats()
}
